# Auto-generated Excel COM-interop script
# Applies a batch of market-data value refreshes to the Goblin_Profits leve-profit tables.
# Source: diff of Sheets/Goblin_Profits.xlsx (scheduled runner data refresh).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 746196.0600000001
$ws.Range("J17").Value = 835584.25
$ws.Range("L17").Value = 2506752.75
$ws.Range("N17").Value = -2507088.75
$ws.Range("H76").Value = 3500
$ws.Range("I76").Value = 5000
$ws.Range("J76").Value = 2000
$ws.Range("K76").Value = 5000
$ws.Range("L76").Value = 2000
$ws.Range("M76").Value = -4685
$ws.Range("N76").Value = -2630
$ws.Range("H79").Value = 3500
$ws.Range("I79").Value = 5000
$ws.Range("J79").Value = 2000
$ws.Range("K79").Value = 5000
$ws.Range("L79").Value = 2000
$ws.Range("M79").Value = -3908
$ws.Range("N79").Value = -4184
$ws.Range("H86").Value = 4966.6665
$ws.Range("I86").Value = 4966.6665
$ws.Range("K86").Value = 4966.6665
$ws.Range("M86").Value = -3843.6665
$ws.Range("H89").Value = 4966.6665
$ws.Range("I89").Value = 4966.6665
$ws.Range("K89").Value = 24833.3325
$ws.Range("M89").Value = -19217.3325
$ws.Range("H103").Value = 1650
$ws.Range("I103").Value = 1350
$ws.Range("J103").Value = 1800
$ws.Range("K103").Value = 4050
$ws.Range("L103").Value = 5400
$ws.Range("M103").Value = -3464
$ws.Range("N103").Value = -6572
$ws.Range("H115").Value = 1310.2222
$ws.Range("I115").Value = 586.5
$ws.Range("K115").Value = 1759.5
$ws.Range("M115").Value = -192.5
$ws.Range("H138").Value = 2229366
$ws.Range("I138").Value = 1779.0952
$ws.Range("K138").Value = 5337.2856
$ws.Range("M138").Value = -197.2856000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1906.4117
$ws.Range("I45").Value = 1357.5
$ws.Range("J45").Value = 2690.5715
$ws.Range("K45").Value = 1357.5
$ws.Range("L45").Value = 2690.5715
$ws.Range("M45").Value = -980.5
$ws.Range("N45").Value = -3444.5715
$ws.Range("H88").Value = 3429.6155
$ws.Range("I88").Value = 3908
$ws.Range("J88").Value = 3130.625
$ws.Range("K88").Value = 3908
$ws.Range("L88").Value = 3130.625
$ws.Range("M88").Value = -3502
$ws.Range("N88").Value = -3942.625
$ws.Range("H91").Value = 3429.6155
$ws.Range("I91").Value = 3908
$ws.Range("J91").Value = 3130.625
$ws.Range("K91").Value = 3908
$ws.Range("L91").Value = 3130.625
$ws.Range("M91").Value = -2504
$ws.Range("N91").Value = -5938.625
$ws.Range("H94").Value = 27000.5
$ws.Range("J94").Value = 27000.5
$ws.Range("L94").Value = 27000.5
$ws.Range("N94").Value = -28802.5
$ws.Range("H110").Value = 655.0714
$ws.Range("I110").Value = 799.7
$ws.Range("J110").Value = 293.5
$ws.Range("K110").Value = 799.7
$ws.Range("L110").Value = 293.5
$ws.Range("M110").Value = 1245.3
$ws.Range("N110").Value = -4383.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2672.775
$ws.Range("J20").Value = 2603.25
$ws.Range("L20").Value = 2603.25
$ws.Range("N20").Value = -3097.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1322.25
$ws.Range("I16").Value = 1096.4615
$ws.Range("K16").Value = 1096.4615
$ws.Range("M16").Value = -809.4614999999999
$ws.Range("H31").Value = 2896.742
$ws.Range("I31").Value = 1626.4286
$ws.Range("K31").Value = 1626.4286
$ws.Range("M31").Value = -1331.4286
$ws.Range("H34").Value = 2896.742
$ws.Range("I34").Value = 1626.4286
$ws.Range("K34").Value = 1626.4286
$ws.Range("M34").Value = -1424.4286
$ws.Range("H86").Value = 9348
$ws.Range("I86").Value = 9685
$ws.Range("K86").Value = 9685
$ws.Range("M86").Value = -8562
$ws.Range("H89").Value = 9348
$ws.Range("I89").Value = 9685
$ws.Range("K89").Value = 48425
$ws.Range("M89").Value = -42809
$ws.Range("H113").Value = 1322.25
$ws.Range("I113").Value = 1096.4615
$ws.Range("K113").Value = 1096.4615
$ws.Range("M113").Value = 1073.5385

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 271.8
$ws.Range("I15").Value = 36.333332
$ws.Range("J15").Value = 625
$ws.Range("K15").Value = 108.999996
$ws.Range("L15").Value = 1875
$ws.Range("M15").Value = 31.000004
$ws.Range("N15").Value = -2155
$ws.Range("H60").Value = 148.33333
$ws.Range("I60").Value = 148.33333
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 444.99999
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -193.99999
$ws.Range("H98").Value = 116
$ws.Range("I98").Value = 116
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 348
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = 1150
$ws.Range("H107").Value = 2098.9285
$ws.Range("J107").Value = 1704.2727
$ws.Range("L107").Value = 5112.8181
$ws.Range("N107").Value = -8952.8181
$ws.Range("H125").Value = 8057.25
$ws.Range("I125").Value = 8057.25
$ws.Range("K125").Value = 24171.75
$ws.Range("M125").Value = -19251.75
$ws.Range("H128").Value = 300000
$ws.Range("I128").Value = 300000
$ws.Range("K128").Value = 900000
$ws.Range("M128").Value = -895020
$ws.Range("H131").Value = 3924882
$ws.Range("I131").Value = 1095.1666
$ws.Range("J131").Value = 6065129.5
$ws.Range("K131").Value = 3285.4998
$ws.Range("L131").Value = 18195388.5
$ws.Range("M131").Value = 1754.5002
$ws.Range("N131").Value = -18205468.5
$ws.Range("H139").Value = 5539.871
$ws.Range("I139").Value = 5056.615
$ws.Range("K139").Value = 15169.845
$ws.Range("M139").Value = -10029.845

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 82.53333000000001
$ws.Range("I2").Value = 41.285713
$ws.Range("J2").Value = 118.625
$ws.Range("K2").Value = 41.285713
$ws.Range("L2").Value = 118.625
$ws.Range("M2").Value = 71.714287
$ws.Range("N2").Value = -344.625
$ws.Range("H80").Value = 47622890
$ws.Range("I80").Value = 100002120
$ws.Range("K80").Value = 100002120
$ws.Range("M80").Value = -100001122
$ws.Range("H83").Value = 47622890
$ws.Range("I83").Value = 100002120
$ws.Range("K83").Value = 500010600
$ws.Range("M83").Value = -500005608
$ws.Range("H97").Value = 1714.6875
$ws.Range("I97").Value = 1233.5
$ws.Range("J97").Value = 2516.6667
$ws.Range("K97").Value = 1233.5
$ws.Range("L97").Value = 2516.6667
$ws.Range("M97").Value = -737.5
$ws.Range("N97").Value = -3508.6667
$ws.Range("H102").Value = 1357.05
$ws.Range("I102").Value = 745.0333000000001
$ws.Range("K102").Value = 745.0333000000001
$ws.Range("M102").Value = 876.9666999999999
$ws.Range("H132").Value = 47621984
$ws.Range("I132").Value = 71431190
$ws.Range("K132").Value = 214293570
$ws.Range("M132").Value = -214291040
$ws.Range("H133").Value = 99998.836
$ws.Range("J133").Value = 99998.836
$ws.Range("L133").Value = 99998.836
$ws.Range("N133").Value = -110118.836

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 29416870
$ws.Range("I40").Value = 4517.952
$ws.Range("J40").Value = 76929130
$ws.Range("K40").Value = 4517.952
$ws.Range("L40").Value = 76929130
$ws.Range("M40").Value = -4381.952
$ws.Range("N40").Value = -76929402
$ws.Range("H61").Value = 3756.0466
$ws.Range("I61").Value = 2771.12
$ws.Range("J61").Value = 5124
$ws.Range("K61").Value = 2771.12
$ws.Range("L61").Value = 5124
$ws.Range("M61").Value = -2569.12
$ws.Range("N61").Value = -5528
$ws.Range("H93").Value = 6847.2354
$ws.Range("I93").Value = 2799.6667
$ws.Range("J93").Value = 7714.5713
$ws.Range("K93").Value = 2799.6667
$ws.Range("L93").Value = 7714.5713
$ws.Range("M93").Value = -1551.6667
$ws.Range("N93").Value = -10210.5713
$ws.Range("H95").Value = 45833
$ws.Range("J95").Value = 45833
$ws.Range("L95").Value = 45833
$ws.Range("N95").Value = -51325
$ws.Range("H97").Value = 50344
$ws.Range("J97").Value = 50344
$ws.Range("L97").Value = 50344
$ws.Range("N97").Value = -52326
$ws.Range("H100").Value = 45459000
$ws.Range("I100").Value = 50002390
$ws.Range("J100").Value = 41672840
$ws.Range("K100").Value = 50002390
$ws.Range("L100").Value = 41672840
$ws.Range("M100").Value = -50001849
$ws.Range("N100").Value = -41673922
$ws.Range("H113").Value = 3756.0466
$ws.Range("I113").Value = 2771.12
$ws.Range("J113").Value = 5124
$ws.Range("K113").Value = 2771.12
$ws.Range("L113").Value = 5124
$ws.Range("M113").Value = -601.1199999999999
$ws.Range("N113").Value = -9464
$ws.Range("H136").Value = 27672.408
$ws.Range("I136").Value = 4197.619
$ws.Range("K136").Value = 12592.857
$ws.Range("M136").Value = -10042.857

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 16298.5
$ws.Range("I58").Value = 9597
$ws.Range("K58").Value = 9597
$ws.Range("M58").Value = -9289
$ws.Range("H95").Value = 33858.8
$ws.Range("J95").Value = 33858.8
$ws.Range("L95").Value = 33858.8
$ws.Range("N95").Value = -39350.8
$ws.Range("H132").Value = 11497004
$ws.Range("I132").Value = 12347708
$ws.Range("K132").Value = 37043124
$ws.Range("M132").Value = -37040594
$ws.Range("H135").Value = 39715
$ws.Range("J135").Value = 39715
$ws.Range("L135").Value = 39715
$ws.Range("N135").Value = -49855
